# Insert a new data row at row 123 (pushing the existing rows 123-227 down to
# 124-228), then populate the new row with the latest "Pepino dulce" price
# observation (week of 2021-09-08 / serial date 44447).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 123..227 down by one row to make room for the new observation.
$ws.Rows.Item(123).Insert()

# Fill in the newly inserted row 123 with the new record's data. The
# surrounding constant columns (A,B,C,E,F,G,H,N,O,Q,R) match every other row
# in this sheet.
$ws.Cells.Item(123, 1).Value = 8
$ws.Cells.Item(123, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(123, 3).Value = "Coquimbo"
$ws.Cells.Item(123, 4).Value = 44447
$ws.Cells.Item(123, 5).Value = 4
$ws.Cells.Item(123, 6).Value = 100112043
$ws.Cells.Item(123, 7).Value = "Pepino dulce"
$ws.Cells.Item(123, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(123, 9).Value = "Especial"
$ws.Cells.Item(123, 10).Value = 600
$ws.Cells.Item(123, 11).Value = 14000
$ws.Cells.Item(123, 12).Value = 15000
$ws.Cells.Item(123, 13).Value = 14500
$ws.Cells.Item(123, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(123, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(123, 16).Value = 806
$ws.Cells.Item(123, 17).Value = 18
$ws.Cells.Item(123, 18).Value = "Hortaliza"
